# t11.3 (Resultado Primario / NFSP - Sergipe) was re-expressed from nominal
# values into real (inflation-adjusted) values: each "Valor" cell in column D
# is multiplied by the deflator for its year (A column); 2022 is the base
# year, so its ratio is 1 and those rows are numerically unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$deflator = @{
    2010 = 2.0257549540190682
    2011 = 1.9020574487400588
    2012 = 1.7971302942988958
    2013 = 1.696835710110893
    2014 = 1.59465843649781
    2015 = 1.44087359258608
    2016 = 1.3556316351529616
    2017 = 1.3168194188501479
    2018 = 1.2692775962090805
    2019 = 1.216877026455523
    2020 = 1.1642811284414059
    2021 = 1.0578509290788947
    2022 = 1.0
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $year = $ws.Cells.Item($r, 1).Value2
    if ($deflator.ContainsKey([int]$year)) {
        $old = $ws.Cells.Item($r, 4).Value2
        $ratio = $deflator[[int]$year]
        $ws.Cells.Item($r, 4).Value = $old * $ratio
    }
}

# The sheet previously carried a stale <sortState> (from an earlier
# Data > Sort on A2:M14); clearing the sort field definitions drops it from
# the saved worksheet XML, matching the cleaned-up file.
$ws.Sort.SortFields.Clear()
